# Atualização de bases das ligas, do dia: 14-05-2024 às 20:19
#
# 1) Several existing rows get their betting-odds data re-shuffled among
#    each other (a data-correction re-sync: ids/teams/odds move between
#    rows that share the same match date, while the row index (col A) and
#    the date (col D) stay put).
# 2) Four brand-new match rows are appended at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Part 1: re-shuffle data among existing rows.
# Mapping: row X's new content (columns B..AB) = row Y's OLD content,
# where Y = $perm[X]. Column A (row index) and column C/D are unaffected
# (C is constant for every data row; D is identical within each group).
# ---------------------------------------------------------------------

$perm = @{
    94  = 97
    95  = 96
    96  = 98
    97  = 95
    98  = 99
    99  = 94
    100 = 101
    101 = 100
    114 = 115
    115 = 114
    173 = 174
    174 = 173
    204 = 206
    205 = 204
    206 = 205
}

$affectedRows = @(94,95,96,97,98,99,100,101,114,115,173,174,204,205,206)

# Snapshot the "before" content of every affected row first, so the
# cycles/swaps below don't clobber data that's still needed as a source.
$snapshot = @{}
foreach ($r in $affectedRows) {
    $snapshot[$r] = $ws.Range("B$r`:AB$r").Value2
}

foreach ($r in $affectedRows) {
    $src = $perm[$r]
    $ws.Range("B$r`:AB$r").Value2 = $snapshot[$src]
}

# ---------------------------------------------------------------------
# Part 2: append four new rows (222-225) at the bottom of the sheet.
# ---------------------------------------------------------------------

# Pick up formatting (styles) from the last existing row and extend it
# down over the new rows, same as Excel does when you continue a table.
$ws.Range("A221:AB221").Copy() | Out-Null
$ws.Range("A222:AB225").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$newRows = @(
    @{ Row=222; A=220; B=8145468; Date=45423.70833333334; Home="Metropolitanos FC"; Away="Deportivo La Guaira"; FTHG=2; FTAG=1; FTR="H";
       Odds=@(1.909,3.2,3.6,2.375,3.2,3,-0.25,1.875,1.925,2.25,1.975,1.825,1.375,-1,-1,0.875,-1,0.9750000000000001,-1) },
    @{ Row=223; A=221; B=8145469; Date=45423.875; Home="UCV"; Away="Inter de Barinas"; FTHG=4; FTAG=0; FTR="H";
       Odds=@(1.833,3,4.333,1.5,3.5,7,-1,1.875,1.925,2,1.775,2.025,0.5,-1,-1,0.875,-1,0.7749999999999999,-1) },
    @{ Row=224; A=222; B=8145479; Date=45424.75; Home="Angostura FC"; Away="Portuguesa"; FTHG=3; FTAG=4; FTR="A";
       Odds=@(2.3,3.1,2.875,2.1,3.1,3.2,-0.25,1.875,1.925,2.25,1.9,1.9,-1,-1,2.2,-1,0.925,0.8999999999999999,-1) },
    @{ Row=225; A=223; B=8145478; Date=45424.85416666666; Home="Carabobo"; Away="Academia Puerto Cabello"; FTHG=1; FTAG=1; FTR="D";
       Odds=@(2.1,3,3.4,2.1,2.875,3.6,-0.25,1.825,1.975,1.75,1.825,1.975,-1,1.875,-1,-0.5,0.4875,0.4125,-0.5) }
)

$oddsCols = @("J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Range("A$r").Value2 = $nr.A
    $ws.Range("B$r").Value2 = $nr.B
    $ws.Range("C$r").Value2 = "Venezuela Primera Division"
    $ws.Range("D$r").Value2 = $nr.Date
    $ws.Range("E$r").Value2 = $nr.Home
    $ws.Range("F$r").Value2 = $nr.Away
    $ws.Range("G$r").Value2 = $nr.FTHG
    $ws.Range("H$r").Value2 = $nr.FTAG
    $ws.Range("I$r").Value2 = $nr.FTR

    for ($i = 0; $i -lt $oddsCols.Length; $i++) {
        $ws.Range($oddsCols[$i] + $r).Value2 = $nr.Odds[$i]
    }
}
